$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Set the Experimental value (row 7, column B) to the text "true".
# A leading apostrophe forces Excel to store this as literal text instead
# of auto-converting it to the Boolean TRUE.
$ws.Range("B7").Value = "'true"

# Update the Date value (row 8, column B) to the new timestamp
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
